$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.096681
$ws.Range("H2").Value = 2.193362
$ws.Range("I2").Value = 0.004066739499659485
$ws.Range("J2").Value = 0.002727650639497791
$ws.Range("M2").Value = 2.993584
$ws.Range("N2").Value = 5.987168
$ws.Range("O2").Value = 0.05975341067179143
$ws.Range("P2").Value = 0.05276159080888739
$ws.Range("Q2").Value = 3.283006694704
$ws.Range("R2").Value = 13.132026778816
$ws.Range("S2").Value = 0.0002430015554183488
$ws.Range("T2").Value = 0.0001439151869107825
$ws.Range("G3").Value = 1.096681
$ws.Range("H3").Value = 2.193362
$ws.Range("I3").Value = 0.004066739499659485
$ws.Range("J3").Value = 0.002727650639497791
$ws.Range("O3").Value = 0.01437707027524159
$ws.Range("P3").Value = 0.01904218747594478
$ws.Range("Q3").Value = 0.7899133695163334
$ws.Range("R3").Value = 4.739480217098
$ws.Range("S3").Value = 0.00005846779957770522
$ws.Range("T3").Value = 0.00005194043484619761
$ws.Range("G4").Value = 1.096681
$ws.Range("H4").Value = 2.193362
$ws.Range("I4").Value = 0.004066739499659485
$ws.Range("J4").Value = 0.002727650639497791
$ws.Range("M4").Value = 9.947454
$ws.Range("N4").Value = 29.842362
$ws.Range("O4").Value = 0.1985560799365424
$ws.Range("P4").Value = 0.2629841842778907
$ws.Range("Q4").Value = 10.909183800174
$ws.Range("R4").Value = 65.45510280104401
$ws.Range("S4").Value = 0.0008074758531754833
$ws.Range("T4").Value = 0.0007173289784233937
$ws.Range("G5").Value = 1.096681
$ws.Range("H5").Value = 2.193362
$ws.Range("I5").Value = 0.004066739499659485
$ws.Range("J5").Value = 0.002727650639497791
$ws.Range("M5").Value = 33.8274285
$ws.Range("N5").Value = 67.65485700000001
$ws.Range("O5").Value = 0.6752121293844308
$ws.Range("P5").Value = 0.596204730060655
$ws.Range("Q5").Value = 37.0978981148085
$ws.Range("R5").Value = 148.391592459234
$ws.Range("S5").Value = 0.002745911837216855
$ws.Range("T5").Value = 0.001626238213221554
$ws.Range("G6").Value = 1.096681
$ws.Range("H6").Value = 2.193362
$ws.Range("I6").Value = 0.004066739499659485
$ws.Range("J6").Value = 0.002727650639497791
$ws.Range("M6").Value = 1.126386666666667
$ws.Range("N6").Value = 3.37916
$ws.Range("O6").Value = 0.02248323249608615
$ws.Range("P6").Value = 0.02977866283320594
$ws.Range("Q6").Value = 1.235286855986667
$ws.Range("R6").Value = 7.41172113592
$ws.Range("S6").Value = 0.00009143344967186129
$ws.Range("T6").Value = 0.00008122578872038328
$ws.Range("G7").Value = 1.096681
$ws.Range("H7").Value = 2.193362
$ws.Range("I7").Value = 0.004066739499659485
$ws.Range("J7").Value = 0.002727650639497791
$ws.Range("M7").Value = 1.483835
$ws.Range("N7").Value = 4.451505
$ws.Range("O7").Value = 0.02961807723590774
$ws.Range("P7").Value = 0.03922864454341624
$ws.Range("Q7").Value = 1.627293651635
$ws.Range("R7").Value = 9.76376190981
$ws.Range("S7").Value = 0.0001204490045992314
$ws.Range("T7").Value = 0.0001070020373754808
$ws.Range("I8").Value = 0.003699735271573896
$ws.Range("J8").Value = 0.003722239381324789
$ws.Range("M8").Value = 2.993584
$ws.Range("N8").Value = 5.987168
$ws.Range("O8").Value = 0.05975341067179143
$ws.Range("P8").Value = 0.05276159080888739
$ws.Range("Q8").Value = 2.986730688362667
$ws.Range("R8").Value = 17.920384130176
$ws.Range("S8").Value = 0.0002210718010592668
$ws.Range("T8").Value = 0.0001963912711301847
$ws.Range("I9").Value = 0.003699735271573896
$ws.Range("J9").Value = 0.003722239381324789
$ws.Range("O9").Value = 0.01437707027524159
$ws.Range("P9").Value = 0.01904218747594478
$ws.Range("S9").Value = 0.00005319135399920791
$ws.Range("T9").Value = 0.00007087958012953135
$ws.Range("I10").Value = 0.003699735271573896
$ws.Range("J10").Value = 0.003722239381324789
$ws.Range("M10").Value = 9.947454
$ws.Range("N10").Value = 29.842362
$ws.Range("O10").Value = 0.1985560799365424
$ws.Range("P10").Value = 0.2629841842778907
$ws.Range("Q10").Value = 9.924680961976001
$ws.Range("R10").Value = 89.32212865778401
$ws.Range("S10").Value = 0.0007346049323266718
$ws.Range("T10").Value = 0.0009788900873847403
$ws.Range("I11").Value = 0.003699735271573896
$ws.Range("J11").Value = 0.003722239381324789
$ws.Range("M11").Value = 33.8274285
$ws.Range("N11").Value = 67.65485700000001
$ws.Range("O11").Value = 0.6752121293844308
$ws.Range("P11").Value = 0.596204730060655
$ws.Range("Q11").Value = 33.749986240354
$ws.Range("R11").Value = 202.499917442124
$ws.Range("S11").Value = 0.002498106130878095
$ws.Range("T11").Value = 0.002219216725563885
$ws.Range("I12").Value = 0.003699735271573896
$ws.Range("J12").Value = 0.003722239381324789
$ws.Range("M12").Value = 1.126386666666667
$ws.Range("N12").Value = 3.37916
$ws.Range("O12").Value = 0.02248323249608615
$ws.Range("P12").Value = 0.02977866283320594
$ws.Range("Q12").Value = 1.123807992124445
$ws.Range("R12").Value = 10.11427192912
$ws.Range("S12").Value = 0.00008318200828476634
$ws.Range("T12").Value = 0.000110843311520952
$ws.Range("I13").Value = 0.003699735271573896
$ws.Range("J13").Value = 0.003722239381324789
$ws.Range("M13").Value = 1.483835
$ws.Range("N13").Value = 4.451505
$ws.Range("O13").Value = 0.02961807723590774
$ws.Range("P13").Value = 0.03922864454341624
$ws.Range("Q13").Value = 1.480438007073333
$ws.Range("R13").Value = 13.32394206366
$ws.Range("S13").Value = 0.0001095790450258878
$ws.Range("T13").Value = 0.0001460184055954957
$ws.Range("G14").Value = 52.66178533333333
$ws.Range("H14").Value = 157.985356
$ws.Range("I14").Value = 0.1952817296348302
$ws.Range("J14").Value = 0.1964695555611368
$ws.Range("M14").Value = 2.993584
$ws.Range("N14").Value = 5.987168
$ws.Range("O14").Value = 0.05975341067179143
$ws.Range("P14").Value = 0.05276159080888739
$ws.Range("Q14").Value = 157.6474779853013
$ws.Range("R14").Value = 945.8848679118079
$ws.Range("S14").Value = 0.01166874938756775
$ws.Range("T14").Value = 0.01036604629692067
$ws.Range("G15").Value = 52.66178533333333
$ws.Range("H15").Value = 157.985356
$ws.Range("I15").Value = 0.1952817296348302
$ws.Range("J15").Value = 0.1964695555611368
$ws.Range("O15").Value = 0.01437707027524159
$ws.Range("P15").Value = 0.01904218747594478
$ws.Range("Q15").Value = 37.93103764668045
$ws.Range("R15").Value = 341.379338820124
$ws.Range("S15").Value = 0.002807579150430681
$ws.Range("T15").Value = 0.003741210110310716
$ws.Range("G16").Value = 52.66178533333333
$ws.Range("H16").Value = 157.985356
$ws.Range("I16").Value = 0.1952817296348302
$ws.Range("J16").Value = 0.1964695555611368
$ws.Range("M16").Value = 9.947454
$ws.Range("N16").Value = 29.842362
$ws.Range("O16").Value = 0.1985560799365424
$ws.Range("P16").Value = 0.2629841842778907
$ws.Range("Q16").Value = 523.850687161208
$ws.Range("R16").Value = 4714.656184450872
$ws.Range("S16").Value = 0.0387743747195196
$ws.Range("T16").Value = 0.05166838580468529
$ws.Range("G17").Value = 52.66178533333333
$ws.Range("H17").Value = 157.985356
$ws.Range("I17").Value = 0.1952817296348302
$ws.Range("J17").Value = 0.1964695555611368
$ws.Range("M17").Value = 33.8274285
$ws.Range("N17").Value = 67.65485700000001
$ws.Range("O17").Value = 0.6752121293844308
$ws.Range("P17").Value = 0.596204730060655
$ws.Range("Q17").Value = 1781.412778045682
$ws.Range("R17").Value = 10688.47666827409
$ws.Range("S17").Value = 0.1318565924966084
$ws.Range("T17").Value = 0.1171360783384644
$ws.Range("G18").Value = 52.66178533333333
$ws.Range("H18").Value = 157.985356
$ws.Range("I18").Value = 0.1952817296348302
$ws.Range("J18").Value = 0.1964695555611368
$ws.Range("M18").Value = 1.126386666666667
$ws.Range("N18").Value = 3.37916
$ws.Range("O18").Value = 0.02248323249608615
$ws.Range("P18").Value = 0.02977866283320594
$ws.Range("Q18").Value = 59.31753284232889
$ws.Range("R18").Value = 533.85779558096
$ws.Range("S18").Value = 0.004390564529617725
$ws.Range("T18").Value = 0.005850600652044914
$ws.Range("G19").Value = 52.66178533333333
$ws.Range("H19").Value = 157.985356
$ws.Range("I19").Value = 0.1952817296348302
$ws.Range("J19").Value = 0.1964695555611368
$ws.Range("M19").Value = 1.483835
$ws.Range("N19").Value = 4.451505
$ws.Range("O19").Value = 0.02961807723590774
$ws.Range("P19").Value = 0.03922864454341624
$ws.Range("Q19").Value = 78.14140024008667
$ws.Range("R19").Value = 703.27260216078
$ws.Range("S19").Value = 0.005783869351086055
$ws.Range("T19").Value = 0.007707234358710803
$ws.Range("G20").Value = 3.794489
$ws.Range("H20").Value = 7.588978
$ws.Range("I20").Value = 0.0140708175826183
$ws.Range("J20").Value = 0.009437603411946896
$ws.Range("M20").Value = 2.993584
$ws.Range("N20").Value = 5.987168
$ws.Range("O20").Value = 0.05975341067179143
$ws.Range("P20").Value = 0.05276159080888739
$ws.Range("Q20").Value = 11.359121558576
$ws.Range("R20").Value = 45.436486234304
$ws.Range("S20").Value = 0.0008407793415020548
$ws.Range("T20").Value = 0.0004979429694377017
$ws.Range("G21").Value = 3.794489
$ws.Range("H21").Value = 7.588978
$ws.Range("I21").Value = 0.0140708175826183
$ws.Range("J21").Value = 0.009437603411946896
$ws.Range("O21").Value = 0.01437707027524159
$ws.Range("P21").Value = 0.01904218747594478
$ws.Range("Q21").Value = 2.733080623793667
$ws.Range("R21").Value = 16.398483742762
$ws.Range("S21").Value = 0.0002022971332154082
$ws.Range("T21").Value = 0.0001797126134939089
$ws.Range("G22").Value = 3.794489
$ws.Range("H22").Value = 7.588978
$ws.Range("I22").Value = 0.0140708175826183
$ws.Range("J22").Value = 0.009437603411946896
$ws.Range("M22").Value = 9.947454
$ws.Range("N22").Value = 29.842362
$ws.Range("O22").Value = 0.1985560799365424
$ws.Range("P22").Value = 0.2629841842778907
$ws.Range("Q22").Value = 37.74550478100601
$ws.Range("R22").Value = 226.473028686036
$ws.Range("S22").Value = 0.002793846380706865
$ws.Range("T22").Value = 0.002481940434829093
$ws.Range("G23").Value = 3.794489
$ws.Range("H23").Value = 7.588978
$ws.Range("I23").Value = 0.0140708175826183
$ws.Range("J23").Value = 0.009437603411946896
$ws.Range("M23").Value = 33.8274285
$ws.Range("N23").Value = 67.65485700000001
$ws.Range("O23").Value = 0.6752121293844308
$ws.Range("P23").Value = 0.596204730060655
$ws.Range("Q23").Value = 128.3578053415365
$ws.Range("R23").Value = 513.431221366146
$ws.Range("S23").Value = 0.00950078670213959
$ws.Range("T23").Value = 0.005626743794639315
$ws.Range("G24").Value = 3.794489
$ws.Range("H24").Value = 7.588978
$ws.Range("I24").Value = 0.0140708175826183
$ws.Range("J24").Value = 0.009437603411946896
$ws.Range("M24").Value = 1.126386666666667
$ws.Range("N24").Value = 3.37916
$ws.Range("O24").Value = 0.02248323249608615
$ws.Range("P24").Value = 0.02977866283320594
$ws.Range("Q24").Value = 4.274061816413333
$ws.Range("R24").Value = 25.64437089848
$ws.Range("S24").Value = 0.0003163574631200242
$ws.Range("T24").Value = 0.0002810392099578806
$ws.Range("G25").Value = 3.794489
$ws.Range("H25").Value = 7.588978
$ws.Range("I25").Value = 0.0140708175826183
$ws.Range("J25").Value = 0.009437603411946896
$ws.Range("M25").Value = 1.483835
$ws.Range("N25").Value = 4.451505
$ws.Range("O25").Value = 0.02961807723590774
$ws.Range("P25").Value = 0.03922864454341624
$ws.Range("Q25").Value = 5.630395585315
$ws.Range("R25").Value = 33.78237351189
$ws.Range("S25").Value = 0.0004167505619343575
$ws.Range("T25").Value = 0.0003702243895889971
$ws.Range("G26").Value = 180.09802
$ws.Range("H26").Value = 540.2940599999999
$ws.Range("I26").Value = 0.6678439142690207
$ws.Range("J26").Value = 0.6719061597109177
$ws.Range("M26").Value = 2.993584
$ws.Range("N26").Value = 5.987168
$ws.Range("O26").Value = 0.05975341067179143
$ws.Range("P26").Value = 0.05276159080888739
$ws.Range("Q26").Value = 539.1385511036799
$ws.Range("R26").Value = 3234.831306622079
$ws.Range("S26").Value = 0.03990595167397346
$ws.Range("T26").Value = 0.03545083786063838
$ws.Range("G27").Value = 180.09802
$ws.Range("H27").Value = 540.2940599999999
$ws.Range("I27").Value = 0.6678439142690207
$ws.Range("J27").Value = 0.6719061597109177
$ws.Range("O27").Value = 0.01437707027524159
$ws.Range("P27").Value = 0.01904218747594478
$ws.Range("Q27").Value = 129.7203414861933
$ws.Range("R27").Value = 1167.48307337574
$ws.Range("S27").Value = 0.009601638888338126
$ws.Range("T27").Value = 0.01279456305945739
$ws.Range("G28").Value = 180.09802
$ws.Range("H28").Value = 540.2940599999999
$ws.Range("I28").Value = 0.6678439142690207
$ws.Range("J28").Value = 0.6719061597109177
$ws.Range("M28").Value = 9.947454
$ws.Range("N28").Value = 29.842362
$ws.Range("O28").Value = 0.1985560799365424
$ws.Range("P28").Value = 0.2629841842778907
$ws.Range("Q28").Value = 1791.51676944108
$ws.Range("R28").Value = 16123.65092496972
$ws.Range("S28").Value = 0.132604469626733
$ws.Range("T28").Value = 0.1767006933228659
$ws.Range("G29").Value = 180.09802
$ws.Range("H29").Value = 540.2940599999999
$ws.Range("I29").Value = 0.6678439142690207
$ws.Range("J29").Value = 0.6719061597109177
$ws.Range("M29").Value = 33.8274285
$ws.Range("N29").Value = 67.65485700000001
$ws.Range("O29").Value = 0.6752121293844308
$ws.Range("P29").Value = 0.596204730060655
$ws.Range("Q29").Value = 6092.252894541571
$ws.Range("R29").Value = 36553.51736724942
$ws.Range("S29").Value = 0.4509363114500187
$ws.Range("T29").Value = 0.400593630576539
$ws.Range("G30").Value = 180.09802
$ws.Range("H30").Value = 540.2940599999999
$ws.Range("I30").Value = 0.6678439142690207
$ws.Range("J30").Value = 0.6719061597109177
$ws.Range("M30").Value = 1.126386666666667
$ws.Range("N30").Value = 3.37916
$ws.Range("O30").Value = 0.02248323249608615
$ws.Range("P30").Value = 0.02977866283320594
$ws.Range("Q30").Value = 202.8600084210667
$ws.Range("R30").Value = 1825.7400757896
$ws.Range("S30").Value = 0.01501528999560662
$ws.Range("T30").Value = 0.02000846698558564
$ws.Range("G31").Value = 180.09802
$ws.Range("H31").Value = 540.2940599999999
$ws.Range("I31").Value = 0.6678439142690207
$ws.Range("J31").Value = 0.6719061597109177
$ws.Range("M31").Value = 1.483835
$ws.Range("N31").Value = 4.451505
$ws.Range("O31").Value = 0.02961807723590774
$ws.Range("P31").Value = 0.03922864454341624
$ws.Range("Q31").Value = 267.2357455067
$ws.Range("R31").Value = 2405.1217095603
$ws.Range("S31").Value = 0.0197802526343508
$ws.Range("T31").Value = 0.02635796790583145
$ws.Range("G32").Value = 31.02214
$ws.Range("H32").Value = 93.06641999999999
$ws.Range("I32").Value = 0.1150370637422974
$ws.Range("J32").Value = 0.1157367912951761
$ws.Range("M32").Value = 2.993584
$ws.Range("N32").Value = 5.987168
$ws.Range("O32").Value = 0.05975341067179143
$ws.Range("P32").Value = 0.05276159080888739
$ws.Range("Q32").Value = 92.86738194975999
$ws.Range("R32").Value = 557.2042916985599
$ws.Range("S32").Value = 0.006873856912270546
$ws.Range("T32").Value = 0.006106457223849681
$ws.Range("G33").Value = 31.02214
$ws.Range("H33").Value = 93.06641999999999
$ws.Range("I33").Value = 0.1150370637422974
$ws.Range("J33").Value = 0.1157367912951761
$ws.Range("O33").Value = 0.01437707027524159
$ws.Range("P33").Value = 0.01904218747594478
$ws.Range("Q33").Value = 22.34451325135333
$ws.Range("R33").Value = 201.10061926218
$ws.Range("S33").Value = 0.001653895949680456
$ws.Range("T33").Value = 0.002203881677707037
$ws.Range("G34").Value = 31.02214
$ws.Range("H34").Value = 93.06641999999999
$ws.Range("I34").Value = 0.1150370637422974
$ws.Range("J34").Value = 0.1157367912951761
$ws.Range("M34").Value = 9.947454
$ws.Range("N34").Value = 29.842362
$ws.Range("O34").Value = 0.1985560799365424
$ws.Range("P34").Value = 0.2629841842778907
$ws.Range("Q34").Value = 308.59131063156
$ws.Range("R34").Value = 2777.32179568404
$ws.Range("S34").Value = 0.02284130842408073
$ws.Range("T34").Value = 0.03043694564970236
$ws.Range("G35").Value = 31.02214
$ws.Range("H35").Value = 93.06641999999999
$ws.Range("I35").Value = 0.1150370637422974
$ws.Range("J35").Value = 0.1157367912951761
$ws.Range("M35").Value = 33.8274285
$ws.Range("N35").Value = 67.65485700000001
$ws.Range("O35").Value = 0.6752121293844308
$ws.Range("P35").Value = 0.596204730060655
$ws.Range("Q35").Value = 1049.39922276699
$ws.Range("R35").Value = 6296.39533660194
$ws.Range("S35").Value = 0.07767442076756914
$ws.Range("T35").Value = 0.06900282241222683
$ws.Range("G36").Value = 31.02214
$ws.Range("H36").Value = 93.06641999999999
$ws.Range("I36").Value = 0.1150370637422974
$ws.Range("J36").Value = 0.1157367912951761
$ws.Range("M36").Value = 1.126386666666667
$ws.Range("N36").Value = 3.37916
$ws.Range("O36").Value = 0.02248323249608615
$ws.Range("P36").Value = 0.02977866283320594
$ws.Range("Q36").Value = 34.94292486746667
$ws.Range("R36").Value = 314.4863238071999
$ws.Range("S36").Value = 0.002586405049785156
$ws.Range("T36").Value = 0.003446486885376172
$ws.Range("G37").Value = 31.02214
$ws.Range("H37").Value = 93.06641999999999
$ws.Range("I37").Value = 0.1150370637422974
$ws.Range("J37").Value = 0.1157367912951761
$ws.Range("M37").Value = 1.483835
$ws.Range("N37").Value = 4.451505
$ws.Range("O37").Value = 0.02961807723590774
$ws.Range("P37").Value = 0.03922864454341624
$ws.Range("Q37").Value = 46.0317371069
$ws.Range("R37").Value = 414.2856339621
$ws.Range("S37").Value = 0.003407176638911407
$ws.Range("T37").Value = 0.004540197446314013
